# IDMA 2026 course info update
# 1. Update "last updated" date text (cached field result) from
#    Tuesday, January 27, 2026 -> Wednesday, January 28, 2026
# 2. Week 7 date "We 12/2" -> "We 11/2"
# 3. Week 13 date "We 26/3" -> "We 25/3"
# (All other hunks in the source diff are purely proofErr/spell-check
#  markup removal and run merges with no visible text change, so no
#  action is required for those.)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Tuesday, January 27, 2026", $false, $false, $false, $false, $false,
    $true, 1, $false, "Wednesday, January 28, 2026", 2)

$d.Content.Find.Execute(
    "We 12/2", $true, $false, $false, $false, $false,
    $true, 1, $false, "We 11/2", 2)

$d.Content.Find.Execute(
    "We 26/3", $true, $false, $false, $false, $false,
    $true, 1, $false, "We 25/3", 2)
